$d = $word.ActiveDocument

# Locate the paragraph that currently ends the "Experiment" section
# (it holds the trailing bookmark "_GoBack" and the closing gramEnd proofErr).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -like "*to compare performance of mesa driver*other platforms.*") {
        $targetIndex = $i
        break
    }
}
$target = $d.Paragraphs($targetIndex)

# Collapse to the very end of that paragraph (right before its paragraph mark,
# i.e. right after "platforms.") and split it into a new paragraph.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# --- "Milestones" (level-0 list item, matches the "Experiment" heading level) ---
$pMilestones = $d.Paragraphs($targetIndex + 1)
$pMilestones.Range.Text = "Milestones"
$pMilestones.Style = "List Paragraph"
$pMilestones.Range.ListFormat.ListLevelNumber = 1
$pMilestones.Range.Font.Size = 14

function Add-ListItem($afterParagraph, $text, $level, $fontSize) {
    $rr = $afterParagraph.Range
    $rr.Collapse(0)
    $rr.InsertParagraphAfter()
    $newPara = $afterParagraph.Next()
    $newPara.Range.Text = $text
    $newPara.Style = "List Paragraph"
    $newPara.Range.ListFormat.ListLevelNumber = $level
    $newPara.Range.Font.Size = $fontSize
    return $newPara
}

$p = Add-ListItem $pMilestones "Port OpenGL demo code from Windows to Ubuntu (VMware version)" 2 12
$p = Add-ListItem $p "Upgrade Ubuntu (VMware guest) graphics driver with latest Mesa solution" 2 12
$p = Add-ListItem $p "Add FPS(frame per second) report" 2 12
$p = Add-ListItem $p "Add logging functionality to OpenGL APIs" 2 12
$p = Add-ListItem $p "Initiate EGL demo code" 2 12
$p = Add-ListItem $p "Add logging functionality to EGL APIs" 2 12
$p = Add-ListItem $p "Set up Ubuntu environment on a bare metal with latest Mesa solution" 2 12
$p = Add-ListItem $p "Log function calls to libDRM (hardware specific interface in user space)" 2 12

# --- Trailing empty paragraph: keeps "List Paragraph" style/indent, drops the
# numbering, and carries the relocated "_GoBack" bookmark. ---
$rr = $p.Range
$rr.Collapse(0)
$rr.InsertParagraphAfter()
$pTail = $p.Next()
$pTail.Style = "List Paragraph"
$pTail.Range.Font.Size = 12
$pTail.Range.ParagraphFormat.LeftIndent = $word.InchesToPoints(1.0)

# Move the "_GoBack" bookmark from the old trailing paragraph to the new one.
$d.Bookmarks("_GoBack").Delete()
$pTail.Range.Bookmarks.Add("_GoBack") | Out-Null

# Refresh proofing marks (proofErr gramStart/gramEnd/spellStart/spellEnd)
# across the document so the new text gets the same treatment Word would
# apply automatically while editing.
$d.CheckGrammar()
$d.CheckSpelling()
